$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.367.94"
$ws.Range("E2").Value = "  +4.31%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.485.26"
$ws.Range("E3").Value = "  +3.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
Set-TextValue "D5" "585.84"
$ws.Range("E5").Value = "  +2.93%  "

# Row 6 - Solana
Set-TextValue "D6" "147.97"
$ws.Range("E6").Value = "  +7.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.62%  "

# Row 9 - Toncoin
Set-TextValue "D9" "7.70"
$ws.Range("E9").Value = "  +0.14%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.56%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.399"
$ws.Range("E11").Value = "  +4.63%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "4.081.53"
$ws.Range("E12").Value = "  +3.70%  "

# Row 13 - Avalanche
Set-TextValue "D13" "29.84"
$ws.Range("E13").Value = "  +7.00%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.11%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.489.92"
$ws.Range("E15").Value = "  +3.59%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +4.06%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "63.373.08"
$ws.Range("E17").Value = "  +4.11%  "

# Row 18 - Polkadot
Set-TextValue "D18" "6.28"
$ws.Range("E18").Value = "  +3.01%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.36"
$ws.Range("E19").Value = "  +6.56%  "

# Row 20 - Uniswap
Set-TextValue "D20" "9.39"
$ws.Range("E20").Value = "  +5.79%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "392.77"
$ws.Range("E21").Value = "  +2.41%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.565"
$ws.Range("E22").Value = "  +3.25%  "

# Row 23 - Litecoin
Set-TextValue "D23" "75.19"
$ws.Range("E23").Value = "  -0.14%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000119"
$ws.Range("E25").Value = "  +8.93%  "

# Row 26 - WrappedeETH
Set-TextValue "D26" "3.630.00"

# Row 27 - Kaspa
Set-TextValue "D27" "0.185"
$ws.Range("E27").Value = "  -2.70%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.82"
$ws.Range("E28").Value = "  +9.89%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.10%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "8.26"
$ws.Range("E30").Value = "  +5.42%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.43%  "

# Row 32 - Fetch.AI
Set-TextValue "D32" "1.44"
$ws.Range("E32").Value = "  +8.05%  "

# Row 33 - USDe
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.83"
$ws.Range("E34").Value = "  +3.84%  "

# Row 35 - EnergySwap
Set-TextValue "D35" "32.48"
$ws.Range("E35").Value = "  +28.48%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "5.34"
$ws.Range("E36").Value = "  +9.01%  "

# Row 37 - Aptos
Set-TextValue "D37" "7.16"
$ws.Range("E37").Value = "  +5.18%  "

# Row 38 - was ImmutableX, now Monero
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D38" "171.20"
$ws.Range("E38").Value = "  +2.22%  "

# Row 39 - was Monero, now ImmutableX
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D39" "1.58"
$ws.Range("E39").Value = "  +9.82%  "

# Row 40 - RenzoRestakedETH
Set-TextValue "D40" "3.523.01"
$ws.Range("E40").Value = "  +3.74%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0768"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  +4.56%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +7.73%  "

# Row 44 - Filecoin
Set-TextValue "D44" "4.50"
$ws.Range("E44").Value = "  +4.08%  "

# Row 45 - OKB
Set-TextValue "D45" "42.43"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  +10.14%  "

# Row 47 - Maker
Set-TextValue "D47" "2.619.11"
$ws.Range("E47").Value = "  +7.24%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "23.89"
$ws.Range("E48").Value = "  +8.31%  "

# Row 49 - dogwifhat
$ws.Range("E49").Value = "  +19.16%  "

# Row 50 - Cosmos
Set-TextValue "D50" "6.76"
$ws.Range("E50").Value = "  +2.24%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +5.54%  "
